$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) store plain text, e.g. "309.76" or
# "  -4.07%  ", not numbers. Cells whose new value would otherwise be
# auto-detected as a number by Excel are explicitly pre-formatted as Text
# ("@") so the literal string (and not a coerced numeric value) is stored.
$ws.Range("D2").Value = "44.161.25"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "2.257.77"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.76"
$ws.Range("E5").Value = "  -4.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.05"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.575"
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.74"
$ws.Range("E10").Value = "  -4.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0824"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.37"
$ws.Range("E12").Value = "  -4.61%  "
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("D14").Value = "2.601.61"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.261.66"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.842"
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("D18").Value = "44.042.20"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("E19").Value = "  -6.12%  "
$ws.Range("D20").Value = "0.0₃0977"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.37"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.53"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.98"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  -6.45%  "
$ws.Range("E25").Value = "  -8.18%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.43"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.21"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.16"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.87"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.56"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0827"
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.108"
$ws.Range("E37").Value = "  -5.18%  "
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.88"
$ws.Range("E39").Value = "  +3.60%  "
$ws.Range("E40").Value = "  -8.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.41"
$ws.Range("E41").Value = "  -10.47%  "
$ws.Range("E42").Value = "  -3.66%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "1.778.12"
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "87.85"
$ws.Range("E45").Value = "  +5.99%  "
$ws.Range("E46").Value = "  -3.59%  "
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.97"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.41"
$ws.Range("E50").Value = "  -5.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.70"
$ws.Range("E51").Value = "  -5.87%  "
